# The "When does the pool open?" question was mis-transcribed during the
# original QA pre-processing pass; correct the wording in the shared
# question/answer table on Sheet1 (row 7, column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "What time does the pool open?"

# Restore the view to the top of the sheet with the edited cell selected.
$ws.Range("A1").Select() | Out-Null
$ws.Range("A7").Select() | Out-Null
